# "lỗi upfile danh sach sinh vien" — the class-roster template exposed a
# "MẬT KHẨU" (password) column that broke the student-list upload, so it is
# removed here. Column B ("MẬT KHẨU") is deleted outright; Excel shifts the
# remaining header cells ("HỌ ", "TÊN") and the trailing blank cell one
# column to the left to fill the gap, and the now-unused shared string is
# dropped automatically because nothing references it any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Columns("B:B").Delete()

# Leave the selection where it was recorded in the saved workbook.
[void]$ws.Range("F12").Select()
